# Applies the "Preliminary postgres code padding data" edit to the
# "alignment test" worksheet (xl/worksheets/sheet2.xml).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("alignment test")

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 0. Seed brand-new shared strings in the exact order the original edit
#    introduced them, by writing the cells that first reference each one
#    before anything else touches the shared-string table. (The cells'
#    styling/structure is fixed up properly in the later sections; this
#    first pass only has to get the *text* in place so string ids line up.)
# ---------------------------------------------------------------------------
$ws.Range("B113").Value2 = "pre_md"
$ws.Range("D113").Value2 = "11.2-release_nopti"
$ws.Range("F113").Value2 = "11.2-release_pti"
$ws.Range("A111").Value2 = "1 phys core, 2 workers, 1G, select-only, linked with lld --no-rosegment and max-page-size=0x200000"
$ws.Range("F118").Value2 = "What does DTLB and STLB look like? The increse here might  be from DTLB competition."
$ws.Range("H68").Value2 = "How many of these are 4K, and how many are 2M?"
$ws.Range("N62").Value2 = "Ms_Lf with padding"
$ws.Range("Q71").Value2 = "ITLB_MISSES.STLB_HIT_4K"
$ws.Range("Q72").Value2 = "ITLB_MISSES.STLB_HIT_2M"

# ---------------------------------------------------------------------------
# 1. Column widths: new columns N (14) and O (15) get explicit widths.
#    (COM ColumnWidth differs from the raw OOXML "width" attribute by a
#    constant offset measured empirically against this workbook's font.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(14).ColumnWidth = 14 - 0.8320313
$ws.Columns.Item(15).ColumnWidth = 19 - 0.8320313

# ---------------------------------------------------------------------------
# 2. Row 61: extend the "LTO+PGO" header band two columns to the right
#    (merged J61:M61 -> J61:O61) and give N61/O61 the same blank header
#    formatting as the rest of the row.
# ---------------------------------------------------------------------------
$ws.Range("J61:M61").UnMerge() | Out-Null
$ws.Range("M61").Copy() | Out-Null
$ws.Range("N61:O61").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("J61:O61").Merge() | Out-Null

# ---------------------------------------------------------------------------
# 3. Row 62: new "Ms_Lf with padding" column group label + merged N62:O62.
# ---------------------------------------------------------------------------
$ws.Range("L62").Copy() | Out-Null
$ws.Range("N62").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("M62").Copy() | Out-Null
$ws.Range("O62").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("N62").Value2 = "Ms_Lf with padding"
$ws.Range("N62:O62").Merge() | Out-Null

# ---------------------------------------------------------------------------
# 4. Row 63: "per transaction" / "total (8000000)" sub-headers for the new
#    column group.
# ---------------------------------------------------------------------------
$ws.Range("L63").Copy() | Out-Null
$ws.Range("N63").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("M63").Copy() | Out-Null
$ws.Range("O63").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("N63").Value2 = "per transaction"
$ws.Range("O63").Value2 = "total (8000000)"

# ---------------------------------------------------------------------------
# 5. Rows 64-67: updated L/M values (re-measured numbers) + new N/O columns
#    (the "Ms_Lf with padding" run) + new Q column (row label echo).
# ---------------------------------------------------------------------------
$ws.Range("L64").Value2 = 138231.978798
$ws.Range("M64").Value2 = 1105855830385
$ws.Range("N64").Value2 = 136578.52948699999
$ws.Range("O64").Value2 = 1092628235895
$ws.Range("Q64").Value2 = "CPU_CLK_UNHALTED.THREAD_P"

$ws.Range("L65").Value2 = 186.79355100000001
$ws.Range("M65").Value2 = 1494348410
$ws.Range("N65").Value2 = 190.30417299999999
$ws.Range("O65").Value2 = 1522433383.5
$ws.Range("Q65").Value2 = "BR_MISP_RETIRED.ALL_BRANCHES"

$ws.Range("L66").Value2 = 50.960827999999999
$ws.Range("M66").Value2 = 407686620.5
$ws.Range("N66").Value2 = 45.757987999999997
$ws.Range("O66").Value2 = 366063906.5
$ws.Range("Q66").Value2 = "ITLB_MISSES.MISS_CAUSES_A_WALK"

$ws.Range("A67").Value2 = "ITLB_MISSES.WALK_COMPLETED"
$ws.Range("L67").Value2 = 21.993981000000002
$ws.Range("M67").Value2 = 175951846
$ws.Range("N67").Value2 = 18.998203
$ws.Range("O67").Value2 = 151985624.5
$ws.Range("Q67").Value2 = "ITLB_MISSES.WALK_COMPLETED"

# ---------------------------------------------------------------------------
# 6. New rows 68-72: follow-up comment + 2M/4K STLB-hit breakdown rows.
# ---------------------------------------------------------------------------
$ws.Range("L69").Value2 = 138358.071375
$ws.Range("M69").Value2 = 1106864571002
$ws.Range("N69").Value2 = 136311.46108099999
$ws.Range("O69").Value2 = 1090491688645
$ws.Range("Q69").Value2 = "CPU_CLK_UNHALTED.THREAD_P"

$ws.Range("L70").Value2 = 21.489221000000001
$ws.Range("M70").Value2 = 171913765
$ws.Range("N70").Value2 = 17.734165000000001
$ws.Range("O70").Value2 = 141873322
$ws.Range("Q70").Value2 = "ITLB_MISSES.WALK_COMPLETED_4K"

$ws.Range("L71").Value2 = 260.47150799999997
$ws.Range("M71").Value2 = 2083772060.5
$ws.Range("N71").Value2 = 51.112361999999997
$ws.Range("O71").Value2 = 408898897.5

$ws.Range("L72").Value2 = 1.9895620000000001
$ws.Range("M72").Value2 = 15916498
$ws.Range("N72").Value2 = 3.3670550000000001
$ws.Range("O72").Value2 = 26936437.5

Write-Output "done part 1"
